# 0.2 version: support nested workflow
#
# - Fix "Excute" -> "Execute" typos in the Messages sheet (Robot1Name,
#   Robot2Name, WhichRobotMessage values).
# - Switch the active/selected tab from Settings to Messages, and update
#   the remembered selection on every sheet (Settings -> A7, Messages ->
#   B9, Log -> B8).

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsMessages = $wb.Worksheets.Item("Messages")
$wsLog      = $wb.Worksheets.Item("Log")

# Fix the typos (order matters so the new shared strings land at the end
# in Robot1Name, Robot2Name, WhichRobotMessage order).
$wsMessages.Range("B9").Value  = "Robot 1: Execute all your test cases"
$wsMessages.Range("B10").Value = "Robot 2: Execute just one test case (you pickup later)"
$wsMessages.Range("B8").Value  = "Please choose which sub-robot you want to execute"

# Update the remembered selection on the sheets that are not becoming the
# active tab.
$wsSettings.Activate()
$wsSettings.Range("A7").Select()

$wsLog.Activate()
$wsLog.Range("B8").Select()

# Messages becomes the active tab, selection on B9.
$wsMessages.Activate()
$wsMessages.Range("B9").Select()
